$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns B, C, E are always non-numeric strings (names, URLs, percentages with spaces),
# so they naturally stay text. Column D sometimes holds values that look numeric
# (e.g. "606.05"); temporarily force Text format so Excel keeps them as strings,
# matching the original inlineStr cell type, then restore the original "Normal" style.
$ws.Range("D2").Value = "64.351.38"
$ws.Range("E2").Value = "  -2.86%  "
$ws.Range("D3").Value = "3.147.99"
$ws.Range("E3").Value = "  -1.71%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "606.05"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.18%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "149.55"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -4.12%  "
$ws.Range("E7").Value = "  +0.13%  "
$ws.Range("D8").Value = "3.147.45"
$ws.Range("E8").Value = "  -1.74%  "
$ws.Range("E9").Value = "  -3.21%  "
$ws.Range("E10").Value = "  -4.51%  "
$ws.Range("E11").Value = "  -1.01%  "
$ws.Range("E12").Value = "  -4.57%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000261"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.69%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.84"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.91%  "
$ws.Range("D15").Value = "3.662.14"
$ws.Range("E15").Value = "  -1.75%  "
$ws.Range("D16").Value = "64.374.25"
$ws.Range("E16").Value = "  -3.03%  "
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "3.154.38"
$ws.Range("E17").Value = "  -1.38%  "
$ws.Range("B18").Value = "TRON"
$ws.Range("C18").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.114"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.10%  "
$ws.Range("E19").Value = "  -3.55%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "483.34"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -4.41%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.714"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.09%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.78"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.68%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.85"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -5.23%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "84.12"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.12%  "
$ws.Range("E26").Value = "  +0.11%  "
$ws.Range("E27").Value = "  -1.73%  "
$ws.Range("E28").Value = "  -5.09%  "
$ws.Range("E29").Value = "  -4.04%  "
$ws.Range("E30").Value = "  -3.04%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.96"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.16%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.73"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -6.55%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "26.78"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.99%  "
$ws.Range("E35").Value = "  -5.47%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.12"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -4.77%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "54.41"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.69%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.26"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +7.07%  "
$ws.Range("D39").Value = "0.0₃0756"
$ws.Range("E39").Value = "  -1.66%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "452.73"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -9.47%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0401"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -4.03%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.123"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -5.80%  "
$ws.Range("E43").Value = "  -2.50%  "
$ws.Range("D44").Value = "2.888.00"
$ws.Range("E44").Value = "  -0.74%  "
$ws.Range("E45").Value = "  -7.58%  "
$ws.Range("E46").Value = "  -4.47%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "26.77"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.86%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.998"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.115"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.92%  "
$ws.Range("E50").Value = "  -2.69%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "119.95"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.80%  "
